$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension -> measure labels for estado-civil and sexo
$ws.Range("B2").Value = "iaest-measure:estado-civil"
$ws.Range("E2").Value = "iaest-measure:sexo"

# Update "dim" -> "medida" for estado-civil and sexo rows (row 3)
$ws.Range("B3").Value = "medida"
$ws.Range("E3").Value = "medida"

# Update "skos:Concept" -> "xsd:int" (row 4)
$ws.Range("B4").Value = "xsd:int"
$ws.Range("E4").Value = "xsd:int"

# Remove row 5 (mapping-estado-civil.xlsx / mapping-sexo.xlsx) entirely
$ws.Rows.Item(5).Delete()
